# AddOrEditUserWindow: wire up the backend call that persists a newly
# introduced clearance level ("Top-secret") to the lookup table on the
# "clearance" sheet, then leave that sheet focused/selected on the new
# entry - mirroring what the window does after a successful save.

$wb = $excel.ActiveWorkbook

$people    = $wb.Worksheets.Item("people")
$clearance = $wb.Worksheets.Item("clearance")

# Append the new clearance level as a third header column.
$clearance.Range("C1").Value = "Top-secret"

# Bring the clearance sheet to the front with the new cell selected,
# just like the UI does right after the add/edit call returns.
$clearance.Activate()
$clearance.Range("C1").Select() | Out-Null
